$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell I1 ("I0"), matching the style of the existing header row (bold, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

# New header cell J1 ("IF"), same style
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# New data cells for row 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
